# ---------------------------------------------------------------------------
# Updates to network creation (parameterisation for genuine islands like Hong
# Kong) and public-open-space identification (boundary=protected_area second
# pass criteria, now requiring an additional 'leisure' tag etc).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsRegion = $wb.Worksheets.Item("region_settings")
$wsOpenSpace = $wb.Worksheets.Item("osm_open_space")
$wsProject = $wb.Worksheets.Item("project_settings")
$wsDestinations = $wb.Worksheets.Item("destinations")

# ---------------------------------------------------------------------------
# 1. region_settings: three new "network" parameter rows (15:17).
#    Shared-string order matters, so values are written in the exact order
#    the original authoring session produced them.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2a. osm_open_space: insert a new column G ("second pass" criteria), pushing
#     the former column G (first pass) to column H, *before* any new text is
#     typed in, so the duplicated values/styles are sourced from the
#     untouched original column F / column G (-> H) content.
# ---------------------------------------------------------------------------

$wsOpenSpace.Columns("G:G").Insert()

# Match formatting of the (now shifted) neighbouring column H.
$wsOpenSpace.Range("H2:H15").Copy()
$wsOpenSpace.Range("G2:G15").PasteSpecial(-4122)

# Copy column F's existing values into the new column G.
$wsOpenSpace.Range("F2:F15").Copy()
$wsOpenSpace.Range("G2:G15").PasteSpecial(-4163)

$wsOpenSpace.Application.CutCopyMode = $false

# "network" (group column) is typed first -> becomes the first brand new
# shared string.
$wsRegion.Range("B15").Value = "network"
$wsRegion.Range("B16").Value = "network"
$wsRegion.Range("B17").Value = "network"

# osm_open_space: new "second pass" column header + value (see section 2),
# these come next in authoring order. F5's original text has already been
# copied into G5 above, so overwriting F5 now only affects the "first pass"
# cell -> "second pass" cell relationship as intended.
$wsOpenSpace.Range("G1").Value = "Second pass criteria (removed boundary=protected_area)"
$wsOpenSpace.Range("F5").Value = "national_park,nature_reserve,forest,state_forest,state_park,regional_park,park,county_park"

# Row 15 : network_not_using_buffered_region
$wsRegion.Range("A15").Value = "network_not_using_buffered_region"
$wsRegion.Range("C15").Value = "Instead of using buffered study region, use regular study region for excerpting network from OSM.  This may allow for looping over true islands to extract individual networks (eg. Hong Kong), which may not be possible with the buffered region (which results in only retaining larget network segment)"
$wsRegion.Range("J15").Value = $true

# Row 17 : network_connection_threshold
$wsRegion.Range("A17").Value = "network_connection_threshold"
$wsRegion.Range("C17").Value = "Minimum distance to retain "
$wsRegion.Range("J17").Value = 200

# Row 16 : network_polygon_iteration
$wsRegion.Range("A16").Value = "network_polygon_iteration"
$wsRegion.Range("C16").Value = "Iterate over polygons for network retrieval, and then combin"
$wsRegion.Range("J16").Value = $true

# ---------------------------------------------------------------------------
# 2. osm_open_space: insert a new column G ("second pass" criteria), pushing
#    the former column G (first pass) to column H. The new column mirrors
#    column F's values (the only exception being row 5, handled above where
#    F5 becomes the new "second pass" text and G5 keeps the former F5 text).
# ---------------------------------------------------------------------------

$wsOpenSpace.Columns("G:G").Insert()

# Match formatting of the (now shifted) neighbouring column H.
$wsOpenSpace.Range("H2:H15").Copy()
$wsOpenSpace.Range("G2:G15").PasteSpecial(-4122)

# Copy column F's existing values into the new column G.
$wsOpenSpace.Range("F2:F15").Copy()
$wsOpenSpace.Range("G2:G15").PasteSpecial(-4163)

# Row 5's F cell is replaced with the new "second pass" text (done above,
# after the value copy so G5 keeps the original F5 text).
$wsOpenSpace.Range("F5").Value = "national_park,nature_reserve,forest,state_forest,state_park,regional_park,park,county_park"

$wsOpenSpace.Rows("1:1").RowHeight = 30

$wsOpenSpace.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. View / selection state.
# ---------------------------------------------------------------------------

$wsProject.Range("D16").Select()

$wsDestinations.Range("B24").Select()

$wsOpenSpace.Range("F1").Select()

$wsRegion.Range("J18").Select()
